# "Adjust move animation and balance"
# Rebalance the "srp" move's startup (B24) and active/hit-advantage (C24)
# frame data on the "Arkusz1" frame-data sheet. Every other changed cell
# (B4, C4, D4, F4, G4, F24, G24 ...) is a formula that recalculates
# automatically from these two raw inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("B24").Value = 13
$ws.Range("C24").Value = 9
